# #5: cash & deposit done
# Update the "存款" (deposit) sheet: turn the F1 header ("局金素梅" label,
# really row1 mirroring row2) into a proper header row, and extend both
# header + data rows with the extra property/legislator metadata columns
# that the other sheets (土地/建物/汽車) already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) -----------------------------------------------
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Copy the existing bold/bordered header formatting (from B1) onto the
# newly added header cells G1:M1 so they visually match B1:F1.
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
# Force the date column to stay a literal text string ("2011-11-22")
# instead of Excel auto-converting it to a date serial value.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2011-11-22"
$ws.Range("I2").NumberFormat = "General"
$ws.Range("J2").Value = "高金素梅"
$ws.Range("K2").Value = 926
$ws.Range("L2").Value = "tmp2f3b1"
$ws.Range("M2").Value = 44

# Copy the existing (non-bold/border-less) data-row formatting (from B2)
# onto the newly added data cells G2:M2 so they match B2:F2.
$ws.Range("B2").Copy()
$ws.Range("G2:M2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
